$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract = @'
The coronavirus disease (COVID-19) caused by the novel severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) has now spread to &gt;200 countries posing a global public health concern.
 Patients with comorbidity, such as hypertension suffer more severe infection with elevated mortality.
 The development of effective antiviral drugs is in urgent need to treat COVID-19 patients.
 Here, we report that calcium channel blockers (CCBs), a type of antihypertensive drug that is widely used in clinics, inhibited the post-entry replication events of SARS-CoV-2 in vitro, while no in vitro anti-SARS-CoV-2 effect was observed for the two other major types of antihypertensive drugs, namely, angiotensin-converting enzyme inhibitors and angiotensin II receptor blockers.
 CCB combined with chloroquine showed a significantly enhanced anti-SARS-CoV-2 efficacy.
 A retrospective clinical investigation on hospitalized COVID-19 patients with hypertension as the only comorbidity revealed that the CCB amlodipine besylate therapy was associated with a decreased case fatality rate.
 The results from this study suggest that CCB administration to COVID-19 patients with hypertension as the comorbidity might improve the disease outcome.

'@

$authors = @'
[Lei-Ke%Zhang%NULL%1,  Yuan%Sun%NULL%2,  Yuan%Sun%NULL%0,  Haolong%Zeng%NULL%2,  Haolong%Zeng%NULL%0,  Qingxing%Wang%NULL%1,  Xiaming%Jiang%NULL%1,  Wei-Juan%Shang%NULL%1,  Yan%Wu%NULL%1,  Shufen%Li%NULL%1,  Yu-Lan%Zhang%NULL%1,  Zhao-Nian%Hao%NULL%1,  Hongbo%Chen%NULL%1,  Runming%Jin%NULL%2,  Runming%Jin%NULL%0,  Wei%Liu%NULL%0,  Hao%Li%lihao_1986@126.com%1,  Ke%Peng%pengke@wh.iov.cn%1,  Gengfu%Xiao%xiaogf@wh.iov.cn%2,  Gengfu%Xiao%xiaogf@wh.iov.cn%0]
'@

$ws.Range("D10").Value = $abstract
$ws.Range("E10").Value = $authors
